$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dois = @(
    "10.1007/s11661-011-0834-5",
    "10.1007/s11665-013-0565-4",
    "10.1007/s11665-023-08603-z",
    "10.1007/s11665-024-10426-5",
    "10.1007/s40430-019-1821-9",
    "10.1016/j.ceramint.2012.08.093",
    "10.1016/j.corsci.2014.07.047",
    "10.1016/j.engfailanal.2006.11.053",
    "10.1016/j.ijfatigue.2017.12.016",
    "10.1016/j.ijrmhm.2005.02.003",
    "10.1016/j.jmrt.2013.03.011",
    "10.1016/j.jmrt.2013.10.005",
    "10.1016/j.jmrt.2015.05.006",
    "10.1016/j.matdes.2008.08.010",
    "10.1016/j.matdes.2015.08.158",
    "10.1016/j.materresbull.2004.04.006",
    "10.1016/j.msea.2003.11.016",
    "10.1016/j.msea.2005.03.073",
    "10.1016/j.msea.2006.05.089",
    "10.1016/j.msea.2007.06.009",
    "10.1016/j.msea.2007.09.024",
    "10.1016/j.msea.2013.04.042",
    "10.1016/j.msea.2016.06.041",
    "10.1016/j.msec.2009.03.006",
    "10.1016/j.proeng.2010.03.178",
    "10.1016/j.surfcoat.2019.04.095",
    "10.1016/j.tafmec.2023.104051",
    "10.1051/matecconf/201816514001",
    "10.1051/matecconf/202032111095",
    "10.1179/mht.2005.053",
    "10.1590/0370-44672019740154",
    "10.1590/1980-5373-MR-2021-0280",
    "10.1590/1980-5373-mr-2015-0737",
    "10.1590/1980-5373-mr-2022-0555",
    "10.1590/S0366-69132014000100006",
    "10.1590/s1516-14392008000100017",
    "10.17563/rbav.v35i1.1017",
    "10.3139/146.111952",
    "10.4028/www.scientific.net/MSF.498-499.104",
    "10.4028/www.scientific.net/MSF.498-499.575",
    "10.4028/www.scientific.net/MSF.530-531.690",
    "10.4028/www.scientific.net/MSF.591-593.30",
    "10.4028/www.scientific.net/MSF.591-593.628",
    "10.4028/www.scientific.net/MSF.591-593.667",
    "10.4028/www.scientific.net/MSF.591-593.839",
    "10.4028/www.scientific.net/MSF.636-637.47",
    "10.4028/www.scientific.net/MSF.636-637.657",
    "10.4028/www.scientific.net/MSF.660-661.1087",
    "10.4028/www.scientific.net/msf.660-661.221",
    "10.4028/www.scientific.net/msf.660-661.235",
    "10.4322/2176-1523.20243029",
    "10.5772/22334"
)

for ($i = 0; $i -lt $dois.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dois[$i]
    $target = "https://doi.org/" + $dois[$i]
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
    $cell.Font.Underline = $true
    $cell.Font.Color = 16711680
}

Write-Output "done"
